$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (OLE color values, BGR order)
$orange = 0x317DED   # ED7D31
$green  = 0x47AD70   # 70AD47
$peach  = 0x84B0F4   # F4B084
$blue   = 0xE7C6B4   # B4C6E7
$yellow = 0x99E6FF   # FFE699
$red    = 0x0000C0   # C00000

$xlRight = -4152

# --- New row labels for 05:30 PM / 06:00 PM (rows 26-27) ---
$ws.Range("B26").Value = "05:30 PM"
$ws.Range("B26").HorizontalAlignment = $xlRight
$ws.Range("B27").Value = "06:00 PM"
$ws.Range("B27").HorizontalAlignment = $xlRight

# --- New orange course block: Tuesday (D) & Friday (G), rows 6-9 ---
$ws.Range("D6").Interior.Color = $orange
$ws.Range("G6").Interior.Color = $orange

$ws.Range("D7").Value = "MOBAPPL"
$ws.Range("D7").Interior.Color = $orange
$ws.Range("G7").Value = "MOBAPPL"
$ws.Range("G7").Interior.Color = $orange

$ws.Range("D8").Value = "07:30 AM - 09:30 AM"
$ws.Range("D8").Interior.Color = $orange
$ws.Range("G8").Value = "07:30 AM - 09:30 AM"
$ws.Range("G8").Interior.Color = $orange

$ws.Range("D9").Interior.Color = $orange
$ws.Range("G9").Interior.Color = $orange

# --- Rename existing green/peach course block (rows 10-13) and extend to row 14 ---
$ws.Range("C11").Value = "DATAMA2"
$ws.Range("F11").Value = "DATAMA2"
$ws.Range("D11").Value = "COMPAIS"
$ws.Range("G11").Value = "COMPAIS"

$ws.Range("C14").Interior.Color = $green
$ws.Range("F14").Interior.Color = $green
$ws.Range("D14").Interior.Color = $peach
$ws.Range("G14").Interior.Color = $peach

# --- Rename existing blue/yellow course block (rows 17-20) and extend to row 21 ---
$ws.Range("C18").Value = "OPESYST"
$ws.Range("F18").Value = "OPESYST"
$ws.Range("D18").Value = "WEBPROG"
$ws.Range("G18").Value = "WEBPROG"

$ws.Range("C21").Interior.Color = $blue
$ws.Range("F21").Interior.Color = $blue
$ws.Range("D21").Interior.Color = $yellow
$ws.Range("G21").Interior.Color = $yellow

# --- Remove old orange/red course block (rows 22-24) entirely ---
$ws.Range("C22:G24").Clear()

# --- New red course block: Thursday (F) only, rows 25-27 ---
$ws.Range("F25").Interior.Color = $red
$ws.Range("F26").Value = "PEDUFOR"
$ws.Range("F26").Interior.Color = $red
$ws.Range("F27").Value = "05:00 PM - 05:40 PM"
$ws.Range("F27").Interior.Color = $red

Write-Output "done"
